$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.069.97"
$ws.Range("E2").Value = "  +3.31%  "
$ws.Range("D3").Value = "1.575.26"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.00%  "
$ws.Range("D5").Value = "'212.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D8").Value = "'23.22"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.41%  "
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").Value = "1.799.15"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").Value = "1.575.40"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").Value = "'0.522"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").Value = "28.050.26"
$ws.Range("E16").Value = "  +3.35%  "
$ws.Range("D17").Value = "'63.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("D18").Value = "'228.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.02%  "
$ws.Range("D19").Value = "0.0₃0707"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").Value = "'9.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "'152.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("D26").Value = "'15.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("D27").Value = "'6.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("D28").Value = "'0.107"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("D34").Value = "1.417.20"
$ws.Range("E34").Value = "  -2.33%  "
$ws.Range("D35").Value = "'1.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E36").Value = "  -4.91%  "
$ws.Range("E37").Value = "  -1.40%  "
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("D39").Value = "'0.542"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("D40").Value = "'2.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.79%  "
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("D43").Value = "'5.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.59%  "
$ws.Range("E44").Value = "  -2.57%  "
$ws.Range("E45").Value = "  +5.00%  "
$ws.Range("D46").Value = "'63.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.43%  "
$ws.Range("D47").Value = "1.712.29"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").Value = "'86.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("E49").Value = "  +2.79%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'38.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.83%  "
